$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("AD2").Value = 6.5
$ws.Range("AG2").Value = 4.75
$ws.Range("AT2").Value = 126
$ws.Range("AX2").Value = 67
$ws.Range("H2").Value = 2.8
$ws.Range("I2").Value = 2.4
$ws.Range("L2").Value = 3.5
$ws.Range("N2").Value = 4.75
$ws.Range("W2").Value = 6

# Row 3
$ws.Range("AJ3").Value = 26
$ws.Range("I3").Value = 2.63
$ws.Range("M3").Value = 1.13
$ws.Range("N3").Value = 6
$ws.Range("Q3").Value = 3.1
$ws.Range("R3").Value = 1.36
$ws.Range("X3").Value = 13

# Row 4
$ws.Range("AC4").Value = 12
$ws.Range("AK4").Value = 41
$ws.Range("AN4").Value = 8.5
$ws.Range("AO4").Value = 19
$ws.Range("AP4").Value = 26
$ws.Range("AR4").Value = 3
$ws.Range("AU4").Value = 6.5
$ws.Range("BA4").Value = 126
$ws.Range("N4").Value = 12
$ws.Range("Q4").Value = 1.8
$ws.Range("R4").Value = 2
$ws.Range("S4").Value = 1.36
$ws.Range("T4").Value = 3
$ws.Range("U4").Value = 1.8
$ws.Range("V4").Value = 1.91
$ws.Range("W4").Value = 7.5
$ws.Range("X4").Value = 8

# Row 5
$ws.Range("AC5").Value = 6.5
$ws.Range("AE5").Value = 19
$ws.Range("AG5").Value = 7
$ws.Range("AI5").Value = 12
$ws.Range("AJ5").Value = 34
$ws.Range("AK5").Value = 29
$ws.Range("AM5").Value = 4.33
$ws.Range("AR5").Value = 2.25
$ws.Range("AV5").Value = 19
$ws.Range("AW5").Value = 34
$ws.Range("AX5").Value = 67
$ws.Range("BA5").Value = 301
$ws.Range("G5").Value = 2.55
$ws.Range("I5").Value = 3
$ws.Range("K5").Value = 1.91
$ws.Range("Q5").Value = 2.6
$ws.Range("R5").Value = 1.48
$ws.Range("S5").Value = 1.57
$ws.Range("T5").Value = 2.25
$ws.Range("W5").Value = 6.5
$ws.Range("Y5").Value = 10
$ws.Range("Z5").Value = 23
